# Update cryptocurrency price (D) and 1h volume change (E) figures
# as refreshed by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.471.26"
$ws.Range("E2").Value = "  -1.73%  "

$ws.Range("D3").Value = "2.635.29"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.38"
$ws.Range("E5").Value = "  -1.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.28"
$ws.Range("E6").Value = "  -2.07%  "

$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.577"
$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").Value = "2.666.83"
$ws.Range("E9").Value = "  +0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.54"
$ws.Range("E10").Value = "  +1.33%  "

$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").Value = "3.095.11"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").Value = "59.309.11"
$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.48"
$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").Value = "2.664.45"
$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.64"
$ws.Range("E19").Value = "  -1.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.78"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.23"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.47"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("E25").Value = "  +0.58%  "

$ws.Range("D26").Value = "2.754.86"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  -2.21%  "

$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0" + [char]0x2083 + "0839"
$ws.Range("E29").Value = "  -0.65%  "

$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.56"
$ws.Range("E32").Value = "  +7.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.18"
$ws.Range("E33").Value = "  -0.55%  "

$ws.Range("E34").Value = "  -2.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.59"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +18.48%  "

$ws.Range("E37").Value = "  +0.99%  "

$ws.Range("E38").Value = "  -1.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.893"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.50"
$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.45"
$ws.Range("E41").Value = "  -0.29%  "

$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "294.94"
$ws.Range("E43").Value = "  -3.34%  "

$ws.Range("E44").Value = "  -0.79%  "

$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("E47").Value = "  -0.93%  "

$ws.Range("E48").Value = "  -1.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.83"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.13"
$ws.Range("E51").Value = "  +0.82%  "
